$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/27/2023  Through  3/5/2023"

# --- Weekly crime statistics table updates (rows 14-30) ---
$ws.Range("C14").Value = 4
$ws.Range("E14").Value = -50
$ws.Range("G14").Value = 38
$ws.Range("H14").Value = -34.210526315789
$ws.Range("I14").Value = 59
$ws.Range("J14").Value = 73
$ws.Range("K14").Value = -19.178082191780
$ws.Range("L14").Value = -7.8125
$ws.Range("M14").Value = -22.368421052631
$ws.Range("N14").Value = -82.544378698224
$ws.Range("C15").Value = 26
$ws.Range("D15").Value = 39
$ws.Range("E15").Value = -33.333333333333
$ws.Range("G15").Value = 138
$ws.Range("H15").Value = -22.463768115942
$ws.Range("I15").Value = 257
$ws.Range("J15").Value = 297
$ws.Range("K15").Value = -13.468013468013
$ws.Range("L15").Value = 17.351598173516
$ws.Range("M15").Value = 27.860696517412
$ws.Range("N15").Value = -51.233396584440
$ws.Range("C16").Value = 269
$ws.Range("D16").Value = 287
$ws.Range("E16").Value = -6.271777003484
$ws.Range("F16").Value = 1165
$ws.Range("G16").Value = 1317
$ws.Range("H16").Value = -11.541381928625
$ws.Range("I16").Value = 2697
$ws.Range("J16").Value = 2752
$ws.Range("K16").Value = -1.998546511627
$ws.Range("L16").Value = 41.798107255520
$ws.Range("M16").Value = -12.576985413290
$ws.Range("N16").Value = -82.041550139832
$ws.Range("C17").Value = 430
$ws.Range("D17").Value = 439
$ws.Range("E17").Value = -2.050113895216
$ws.Range("F17").Value = 1853
$ws.Range("G17").Value = 1758
$ws.Range("H17").Value = 5.403868031854
$ws.Range("I17").Value = 4238
$ws.Range("J17").Value = 3837
$ws.Range("K17").Value = 10.450873077925
$ws.Range("L17").Value = 34.071496361910
$ws.Range("M17").Value = 65.611567018366
$ws.Range("N17").Value = -29.624709398870
$ws.Range("C18").Value = 259
$ws.Range("D18").Value = 280
$ws.Range("E18").Value = -7.5
$ws.Range("F18").Value = 1081
$ws.Range("G18").Value = 1224
$ws.Range("H18").Value = -11.683006535947
$ws.Range("I18").Value = 2545
$ws.Range("J18").Value = 2638
$ws.Range("K18").Value = -3.525398028809
$ws.Range("L18").Value = 20.558976788252
$ws.Range("M18").Value = -20.543240711832
$ws.Range("N18").Value = -85.506008314824
$ws.Range("C19").Value = 933
$ws.Range("D19").Value = 926
$ws.Range("E19").Value = 0.755939524838
$ws.Range("F19").Value = 3598
$ws.Range("G19").Value = 3807
$ws.Range("H19").Value = -5.489887050170
$ws.Range("I19").Value = 8238
$ws.Range("J19").Value = 8569
$ws.Range("K19").Value = -3.862761115649
$ws.Range("L19").Value = 59.434875169343
$ws.Range("M19").Value = 38.686868686868
$ws.Range("N19").Value = -38.444295001120
$ws.Range("C20").Value = 267
$ws.Range("D20").Value = 235
$ws.Range("E20").Value = 13.617021276595
$ws.Range("F20").Value = 1129
$ws.Range("G20").Value = 1072
$ws.Range("H20").Value = 5.317164179104
$ws.Range("I20").Value = 2527
$ws.Range("J20").Value = 2405
$ws.Range("K20").Value = 5.072765072765
$ws.Range("L20").Value = 103.298471440064
$ws.Range("M20").Value = 57.151741293532
$ws.Range("N20").Value = -87.486382093691
$ws.Range("C21").Value = 2188
$ws.Range("D21").Value = 2214
$ws.Range("E21").Value = -1.174345076784
$ws.Range("F21").Value = 8958
$ws.Range("G21").Value = 9354
$ws.Range("H21").Value = -4.233483001924
$ws.Range("I21").Value = 20561
$ws.Range("J21").Value = 20571
$ws.Range("K21").Value = -0.048612123863
$ws.Range("L21").Value = 48.272878055816
$ws.Range("M21").Value = 23.326535508637
$ws.Range("N21").Value = -71.850056817403
$ws.Range("C22").Value = 30
$ws.Range("D22").Value = 55
$ws.Range("E22").Value = -45.454545454545
$ws.Range("F22").Value = 164
$ws.Range("G22").Value = 200
$ws.Range("H22").Value = -18
$ws.Range("I22").Value = 336
$ws.Range("J22").Value = 428
$ws.Range("K22").Value = -21.495327102803
$ws.Range("L22").Value = 42.978723404255
$ws.Range("M22").Value = -6.925207756232
$ws.Range("C23").Value = 126
$ws.Range("D23").Value = 98
$ws.Range("E23").Value = 28.571428571428
$ws.Range("F23").Value = 457
$ws.Range("G23").Value = 445
$ws.Range("H23").Value = 2.696629213483
$ws.Range("I23").Value = 1023
$ws.Range("J23").Value = 965
$ws.Range("K23").Value = 6.010362694300
$ws.Range("L23").Value = 24
$ws.Range("M23").Value = 62.380952380952
$ws.Range("C24").Value = 2008
$ws.Range("D24").Value = 2116
$ws.Range("E24").Value = -5.103969754253
$ws.Range("F24").Value = 8092
$ws.Range("G24").Value = 8417
$ws.Range("H24").Value = -3.861233218486
$ws.Range("I24").Value = 18442
$ws.Range("J24").Value = 17569
$ws.Range("K24").Value = 4.968979452444
$ws.Range("L24").Value = 39.205917874396
$ws.Range("M24").Value = 47.784277586345
$ws.Range("C25").Value = 760
$ws.Range("D25").Value = 811
$ws.Range("E25").Value = -6.288532675709
$ws.Range("F25").Value = 2997
$ws.Range("G25").Value = 2983
$ws.Range("H25").Value = 0.469326181696
$ws.Range("I25").Value = 6852
$ws.Range("J25").Value = 6404
$ws.Range("K25").Value = 6.995627732667
$ws.Range("L25").Value = 39.098660170523
$ws.Range("M25").Value = -2.365346252493
$ws.Range("C26").Value = 45
$ws.Range("D26").Value = 69
$ws.Range("E26").Value = -34.782608695652
$ws.Range("F26").Value = 171
$ws.Range("G26").Value = 231
$ws.Range("H26").Value = -25.974025974026
$ws.Range("I26").Value = 416
$ws.Range("J26").Value = 466
$ws.Range("K26").Value = -10.729613733905
$ws.Range("L26").Value = 12.737127371273
$ws.Range("C27").Value = 95
$ws.Range("D27").Value = 84
$ws.Range("E27").Value = 13.095238095238
$ws.Range("F27").Value = 365
$ws.Range("G27").Value = 369
$ws.Range("H27").Value = -1.084010840108
$ws.Range("I27").Value = 830
$ws.Range("J27").Value = 744
$ws.Range("K27").Value = 11.559139784946
$ws.Range("L27").Value = 30.503144654088
$ws.Range("C28").Value = 13
$ws.Range("D28").Value = 21
$ws.Range("E28").Value = -38.095238095238
$ws.Range("F28").Value = 82
$ws.Range("G28").Value = 81
$ws.Range("H28").Value = 1.234567901234
$ws.Range("I28").Value = 179
$ws.Range("J28").Value = 208
$ws.Range("K28").Value = -13.942307692307
$ws.Range("L28").Value = -0.555555555555
$ws.Range("M28").Value = -14.354066985645
$ws.Range("N28").Value = -80.814576634512
$ws.Range("C29").Value = 13
$ws.Range("D29").Value = 19
$ws.Range("E29").Value = -31.578947368421
$ws.Range("F29").Value = 65
$ws.Range("G29").Value = 73
$ws.Range("H29").Value = -10.958904109589
$ws.Range("I29").Value = 148
$ws.Range("J29").Value = 187
$ws.Range("K29").Value = -20.855614973262
$ws.Range("L29").Value = -11.377245508982
$ws.Range("M29").Value = -15.909090909090
$ws.Range("N29").Value = -82.850521436848
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 24
$ws.Range("E30").Value = -91.666666666666
$ws.Range("F30").Value = 24
$ws.Range("G30").Value = 96
$ws.Range("H30").Value = -75
$ws.Range("I30").Value = 60
$ws.Range("J30").Value = 142
$ws.Range("K30").Value = -57.746478873239
$ws.Range("L30").Value = 50
